$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1, matching the style of the existing header row (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Per-row "time_taken" timestamps for data rows 2-35
$timestamps = @(
    "2021-10-05 13:39:13.562989",
    "2021-10-05 13:39:13.563002",
    "2021-10-05 13:39:13.563006",
    "2021-10-05 13:39:13.563009",
    "2021-10-05 13:39:13.563012",
    "2021-10-05 13:39:13.563015",
    "2021-10-05 13:39:13.563018",
    "2021-10-05 13:39:13.563021",
    "2021-10-05 13:39:13.563024",
    "2021-10-05 13:39:13.563027",
    "2021-10-05 13:39:13.563030",
    "2021-10-05 13:39:13.563033",
    "2021-10-05 13:39:13.563036",
    "2021-10-05 13:39:13.563039",
    "2021-10-05 13:39:13.563042",
    "2021-10-05 13:39:13.563045",
    "2021-10-05 13:39:13.563049",
    "2021-10-05 13:39:13.563052",
    "2021-10-05 13:39:13.563055",
    "2021-10-05 13:39:13.563058",
    "2021-10-05 13:39:13.563061",
    "2021-10-05 13:39:13.563064",
    "2021-10-05 13:39:13.563067",
    "2021-10-05 13:39:13.563070",
    "2021-10-05 13:39:13.563074",
    "2021-10-05 13:39:13.563077",
    "2021-10-05 13:39:13.563080",
    "2021-10-05 13:39:13.563083",
    "2021-10-05 13:39:13.563086",
    "2021-10-05 13:39:13.563089",
    "2021-10-05 13:39:13.563092",
    "2021-10-05 13:39:13.563095",
    "2021-10-05 13:39:13.563099",
    "2021-10-05 13:39:13.563102"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
